$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new DAMSLTag (col I) / new DialogAct (col J)
$changes = @{
    15  = @("sd", "Statement-non-opinion")
    17  = @("aa", "Agree/Accept")
    20  = @("sd", "Statement-non-opinion")
    22  = @("aa", "Agree/Accept")
    28  = @("sv", "Statement-opinion")
    38  = @("sv", "Statement-opinion")
    39  = @("sv", "Statement-opinion")
    57  = @("sv", "Statement-opinion")
    75  = @("aa", "Agree/Accept")
    78  = @("aa", "Agree/Accept")
    79  = @("aa", "Agree/Accept")
    85  = @("%", "Uninterpretable")
    97  = @("sv", "Statement-opinion")
    132 = @("b", "Acknowledge (Backchannel)")
    143 = @("sv", "Statement-opinion")
    144 = @("sv", "Statement-opinion")
    148 = @("sd", "Statement-non-opinion")
    151 = @("aa", "Agree/Accept")
    152 = @("%", "Uninterpretable")
    156 = @("sd", "Statement-non-opinion")
    160 = @("sd", "Statement-non-opinion")
    166 = @("sd", "Statement-non-opinion")
    188 = @("aa", "Agree/Accept")
    197 = @("sd", "Statement-non-opinion")
}

foreach ($row in $changes.Keys) {
    $values = $changes[$row]
    $ws.Range("I$row").Value = $values[0]
    $ws.Range("J$row").Value = $values[1]
}

$wb.Save()
